$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 1.2
$ws.Range("K2").Value = 4.33
$ws.Range("AE2").Value = 13
$ws.Range("J3").Value = 1.17
$ws.Range("K3").Value = 5
$ws.Range("L3").Value = 1.67
$ws.Range("M3").Value = 2.1
$ws.Range("U3").Value = 8
$ws.Range("Z3").Value = 5
$ws.Range("G4").Value = 2.52
$ws.Range("H4").Value = 3.75
$ws.Range("I4").Value = 2.45
$ws.Range("T4").Value = 10
$ws.Range("U4").Value = 15
$ws.Range("V4").Value = 10.25
$ws.Range("W4").Value = 29
$ws.Range("X4").Value = 20
$ws.Range("AD4").Value = 10
$ws.Range("AE4").Value = 14.5
$ws.Range("AF4").Value = 10.25
$ws.Range("AG4").Value = 28
$ws.Range("AH4").Value = 19.5
$ws.Range("G5").Value = 2.9
$ws.Range("H5").Value = 3.05
$ws.Range("I5").Value = 2.52
$ws.Range("J5").Value = 1.09
$ws.Range("K5").Value = 6.5
$ws.Range("L5").Value = 1.4
$ws.Range("M5").Value = 2.8
$ws.Range("N5").Value = 2.2
$ws.Range("O5").Value = 1.62
$ws.Range("P5").Value = 1.47
$ws.Range("Q5").Value = 2.55
$ws.Range("R5").Value = 1.85
$ws.Range("S5").Value = 1.85
$ws.Range("T5").Value = 7.9
$ws.Range("U5").Value = 15
$ws.Range("X5").Value = 29
$ws.Range("Y5").Value = 40
$ws.Range("Z5").Value = 6.5
$ws.Range("AA5").Value = 6.2
$ws.Range("AB5").Value = 16
$ws.Range("AC5").Value = 90
$ws.Range("AD5").Value = 7.1
$ws.Range("AE5").Value = 12.5
$ws.Range("AF5").Value = 10.25
$ws.Range("AG5").Value = 30
$ws.Range("AH5").Value = 25
$ws.Range("AI5").Value = 40
$ws.Range("AJ5").Value = 800
$ws.Range("J6").Value = 1.14
$ws.Range("K6").Value = 5.5
$ws.Range("L6").Value = 1.62
$ws.Range("M6").Value = 2.2
$ws.Range("N6").Value = 3.1
$ws.Range("O6").Value = 1.36
$ws.Range("G8").Value = 1.95
$ws.Range("H8").Value = 3.3
$ws.Range("I8").Value = 3.65
$ws.Range("L8").Value = 1.28
$ws.Range("M8").Value = 3
$ws.Range("N8").Value = 1.85
$ws.Range("O8").Value = 1.75
$ws.Range("P8").Value = 1.4
$ws.Range("Q8").Value = 2.52
$ws.Range("R8").Value = 1.7
$ws.Range("S8").Value = 1.91
$ws.Range("T8").Value = 7.1
$ws.Range("U8").Value = 9.25
$ws.Range("V8").Value = 8.5
$ws.Range("X8").Value = 16
$ws.Range("Y8").Value = 27
$ws.Range("Z8").Value = 9.75
$ws.Range("AA8").Value = 6.4
$ws.Range("AB8").Value = 14
$ws.Range("AC8").Value = 65
$ws.Range("AE8").Value = 21
$ws.Range("AJ8").Value = 500
$ws.Range("T9").Value = 6.7
$ws.Range("U9").Value = 10.75
$ws.Range("X9").Value = 21
$ws.Range("AD9").Value = 7.5
$ws.Range("AE9").Value = 14
$ws.Range("AF9").Value = 11.25
$ws.Range("AG9").Value = 40
$ws.Range("G10").Value = 2.2
$ws.Range("I10").Value = 3.3
$ws.Range("R10").Value = 2.25
$ws.Range("S10").Value = 1.57
$ws.Range("X10").Value = 23
$ws.Range("AD10").Value = 7
$ws.Range("M11").Value = 3.75
$ws.Range("T11").Value = 8.5
$ws.Range("U11").Value = 8.75
$ws.Range("W11").Value = 13
$ws.Range("X11").Value = 12
$ws.Range("Y11").Value = 21
$ws.Range("AA11").Value = 7.8
$ws.Range("AD11").Value = 16
$ws.Range("AE11").Value = 30
$ws.Range("AI11").Value = 37
$ws.Range("T12").Value = 13
$ws.Range("U12").Value = 20
$ws.Range("Y12").Value = 27
$ws.Range("Z12").Value = 15
$ws.Range("AB12").Value = 12
$ws.Range("AD12").Value = 10.25
$ws.Range("AE12").Value = 11.75
$ws.Range("AG12").Value = 19
$ws.Range("AI12").Value = 19.5
$ws.Range("G13").Value = 2.25
$ws.Range("I13").Value = 2.92
$ws.Range("W13").Value = 23
$ws.Range("AI13").Value = 27
$ws.Range("G19").Value = 2.75
$ws.Range("H19").Value = 3.5
$ws.Range("I19").Value = 2.35
$ws.Range("J19").Value = 1.04
$ws.Range("K19").Value = 13
$ws.Range("N19").Value = 1.8
$ws.Range("O19").Value = 2
$ws.Range("P19").Value = 1.33
$ws.Range("Q19").Value = 3.25
$ws.Range("R19").Value = 1.62
$ws.Range("S19").Value = 2.2
$ws.Range("T19").Value = 11
$ws.Range("X19").Value = 21
$ws.Range("Y19").Value = 26
$ws.Range("Z19").Value = 13
$ws.Range("AA19").Value = 7
$ws.Range("AD19").Value = 9.5
$ws.Range("AE19").Value = 13
$ws.Range("AI19").Value = 23
$ws.Range("AJ19").Value = 151
$ws.Range("N20").Value = 2.35
$ws.Range("O20").Value = 1.57
$ws.Range("R20").Value = 2.1
$ws.Range("S20").Value = 1.67
$ws.Range("W20").Value = 17
$ws.Range("X20").Value = 19
$ws.Range("AB20").Value = 19
$ws.Range("AD20").Value = 9.5
$ws.Range("AE20").Value = 19
$ws.Range("AJ20").Value = 1250
$ws.Range("I21").Value = 3.65
$ws.Range("T21").Value = 5.8
$ws.Range("U21").Value = 9
$ws.Range("V21").Value = 9.25
$ws.Range("X21").Value = 21
$ws.Range("AB21").Value = 17
$ws.Range("AD21").Value = 8.25
$ws.Range("AE21").Value = 18
$ws.Range("K22").Value = 6.7
$ws.Range("U22").Value = 7.6
$ws.Range("V22").Value = 8
$ws.Range("Y22").Value = 28
$ws.Range("Z22").Value = 6.7
$ws.Range("AB22").Value = 15.5
$ws.Range("AD22").Value = 13
$ws.Range("AE22").Value = 32
$ws.Range("AH22").Value = 55
$ws.Range("AI22").Value = 55
$ws.Range("G30").Value = 2.9
$ws.Range("I30").Value = 2.25
$ws.Range("U30").Value = 17
$ws.Range("V30").Value = 11
$ws.Range("AB30").Value = 13
$ws.Range("AD30").Value = 9.5
$ws.Range("AE30").Value = 12
$ws.Range("AF30").Value = 9
$ws.Range("AG30").Value = 21
$ws.Range("H31").Value = 3.5
$ws.Range("I31").Value = 3.25
$ws.Range("AA31").Value = 6.5
$ws.Range("AD31").Value = 10
$ws.Range("AG31").Value = 34
$ws.Range("H32").Value = 3.35
$ws.Range("Q32").Value = 2.85
$ws.Range("S32").Value = 2.18
$ws.Range("T32").Value = 10
$ws.Range("AD32").Value = 9.5
$ws.Range("H35").Value = 3.6
$ws.Range("I35").Value = 5
$ws.Range("J35").Value = 1.06
$ws.Range("K35").Value = 7.5
$ws.Range("L35").Value = 1.28
$ws.Range("M35").Value = 3.35
$ws.Range("N35").Value = 1.82
$ws.Range("O35").Value = 1.88
$ws.Range("P35").Value = 1.4
$ws.Range("Q35").Value = 2.72
$ws.Range("R35").Value = 1.82
$ws.Range("T35").Value = 6.9
$ws.Range("U35").Value = 7.7
$ws.Range("V35").Value = 7.9
$ws.Range("Z35").Value = 7.5
$ws.Range("AA35").Value = 7.1
$ws.Range("AE35").Value = 30
$ws.Range("AG35").Value = 100
$ws.Range("AH35").Value = 55
$ws.Range("I36").Value = 3.15
$ws.Range("L36").Value = 1.31
$ws.Range("U36").Value = 10.25
$ws.Range("W36").Value = 20
$ws.Range("AD36").Value = 9.5
$ws.Range("AG36").Value = 40
$ws.Range("H38").Value = 3.3
$ws.Range("L38").Value = 1.3
$ws.Range("G39").Value = 3.05
$ws.Range("H39").Value = 3.4
$ws.Range("W39").Value = 37
$ws.Range("AE39").Value = 10
$ws.Range("AG39").Value = 19.5
$ws.Range("G41").Value = 1.52
$ws.Range("I41").Value = 5.9
$ws.Range("L41").Value = 1.24
$ws.Range("M41").Value = 3.3
$ws.Range("N41").Value = 1.72
$ws.Range("O41").Value = 1.9
$ws.Range("S41").Value = 1.82
$ws.Range("T41").Value = 6.9
$ws.Range("W41").Value = 11
$ws.Range("X41").Value = 12
$ws.Range("Y41").Value = 25
$ws.Range("Z41").Value = 11
$ws.Range("AA41").Value = 7.6
$ws.Range("AD41").Value = 16
$ws.Range("AE41").Value = 37
$ws.Range("G42").Value = 2.32
$ws.Range("H42").Value = 3.25
$ws.Range("I42").Value = 2.87
$ws.Range("L42").Value = 1.38
$ws.Range("M42").Value = 2.6
$ws.Range("N42").Value = 2.1
$ws.Range("O42").Value = 1.57
$ws.Range("R42").Value = 1.87
$ws.Range("S42").Value = 1.75
$ws.Range("T42").Value = 6.8
$ws.Range("U42").Value = 10.25
$ws.Range("V42").Value = 9.5
$ws.Range("W42").Value = 23
$ws.Range("X42").Value = 21
$ws.Range("Y42").Value = 37
$ws.Range("AA42").Value = 6.3
$ws.Range("AD42").Value = 7.8
$ws.Range("AE42").Value = 13.5
$ws.Range("AF42").Value = 11
$ws.Range("AG42").Value = 35
$ws.Range("AH42").Value = 27
$ws.Range("AI42").Value = 40
$ws.Range("AJ42").Value = 900
